$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.409.23"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "1.695.59"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.40"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5481"
$ws.Range("E6").Value = "  +4.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.010"
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2737"
$ws.Range("E8").Value = "  +1.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06443"
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.97"
$ws.Range("E10").Value = "  -0.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07678"
$ws.Range("E11").Value = "  +2.17%  "
$ws.Range("D12").Value = "1.696.16"
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5843"
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008415"
$ws.Range("E15").Value = "  -0.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.63"
$ws.Range("E16").Value = "  +2.02%  "
$ws.Range("D17").Value = "26.444.20"
$ws.Range("E17").Value = "  +0.47%  "
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("E20").Value = "  +1.19%  "
$ws.Range("E21").Value = "  +0.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.260"
$ws.Range("E22").Value = "  +0.86%  "
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "148.76"
$ws.Range("E24").Value = "  +2.53%  "
$ws.Range("E25").Value = "  +5.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.917"
$ws.Range("E26").Value = "  +2.24%  "
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06217"
$ws.Range("E28").Value = "  -6.14%  "
$ws.Range("E29").Value = "  +2.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.329"
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("E31").Value = "  +1.06%  "
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("E33").Value = "  +1.36%  "
$ws.Range("E34").Value = "  +1.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6151"
$ws.Range("E35").Value = "  -0.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.413"
$ws.Range("E36").Value = "  +0.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.762"
$ws.Range("E37").Value = "  +1.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01655"
$ws.Range("E38").Value = "  +2.04%  "
$ws.Range("D39").Value = "1.118.30"
$ws.Range("E39").Value = "  +1.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.104"
$ws.Range("E40").Value = "  -4.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8816"
$ws.Range("E41").Value = "  +0.62%  "
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.24"
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("D44").Value = "1.847.51"
$ws.Range("E44").Value = "  +1.07%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000109"
$ws.Range("E45").Value = "  -2.11%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.62"
$ws.Range("E46").Value = "  +1.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.213"
$ws.Range("E47").Value = "  +0.68%  "
$ws.Range("E48").Value = "  -0.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05282"
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.124"
$ws.Range("E50").Value = "  +1.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4306"
$ws.Range("E51").Value = "  +0.12%  "
